# Defects of innate immunity - add "metadata" tab with panel metadata,
# and refresh the "time_taken" timestamps on the "data" tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

$ws1.Range("F2").Value = "2021-10-05 14:33:41.419040"
$ws1.Range("F3").Value = "2021-10-05 14:33:41.419048"
$ws1.Range("F4").Value = "2021-10-05 14:33:41.419051"
$ws1.Range("F5").Value = "2021-10-05 14:33:41.419054"
$ws1.Range("F6").Value = "2021-10-05 14:33:41.419057"
$ws1.Range("F7").Value = "2021-10-05 14:33:41.419060"
$ws1.Range("F8").Value = "2021-10-05 14:33:41.419062"
$ws1.Range("F9").Value = "2021-10-05 14:33:41.419065"
$ws1.Range("F10").Value = "2021-10-05 14:33:41.419068"
$ws1.Range("F11").Value = "2021-10-05 14:33:41.419070"
$ws1.Range("F12").Value = "2021-10-05 14:33:41.419073"
$ws1.Range("F13").Value = "2021-10-05 14:33:41.419076"
$ws1.Range("F14").Value = "2021-10-05 14:33:41.419078"
$ws1.Range("F15").Value = "2021-10-05 14:33:41.419081"
$ws1.Range("F16").Value = "2021-10-05 14:33:41.419083"
$ws1.Range("F17").Value = "2021-10-05 14:33:41.419086"
$ws1.Range("F18").Value = "2021-10-05 14:33:41.419089"
$ws1.Range("F19").Value = "2021-10-05 14:33:41.419092"
$ws1.Range("F20").Value = "2021-10-05 14:33:41.419095"
$ws1.Range("F21").Value = "2021-10-05 14:33:41.419098"
$ws1.Range("F22").Value = "2021-10-05 14:33:41.419100"
$ws1.Range("F23").Value = "2021-10-05 14:33:41.419103"
$ws1.Range("F24").Value = "2021-10-05 14:33:41.419105"
$ws1.Range("F25").Value = "2021-10-05 14:33:41.419108"
$ws1.Range("F26").Value = "2021-10-05 14:33:41.419111"
$ws1.Range("F27").Value = "2021-10-05 14:33:41.419114"
$ws1.Range("F28").Value = "2021-10-05 14:33:41.419117"
$ws1.Range("F29").Value = "2021-10-05 14:33:41.419119"
$ws1.Range("F30").Value = "2021-10-05 14:33:41.419122"
$ws1.Range("F31").Value = "2021-10-05 14:33:41.419124"
$ws1.Range("F32").Value = "2021-10-05 14:33:41.419127"
$ws1.Range("F33").Value = "2021-10-05 14:33:41.419130"
$ws1.Range("F34").Value = "2021-10-05 14:33:41.419133"
$ws1.Range("F35").Value = "2021-10-05 14:33:41.419135"
$ws1.Range("F36").Value = "2021-10-05 14:33:41.419138"
$ws1.Range("F37").Value = "2021-10-05 14:33:41.419140"
$ws1.Range("F38").Value = "2021-10-05 14:33:41.419143"
$ws1.Range("F39").Value = "2021-10-05 14:33:41.419145"
$ws1.Range("F40").Value = "2021-10-05 14:33:41.419148"
$ws1.Range("F41").Value = "2021-10-05 14:33:41.419150"
$ws1.Range("F42").Value = "2021-10-05 14:33:41.419154"
$ws1.Range("F43").Value = "2021-10-05 14:33:41.419156"
$ws1.Range("F44").Value = "2021-10-05 14:33:41.419159"
$ws1.Range("F45").Value = "2021-10-05 14:33:41.419162"
$ws1.Range("F46").Value = "2021-10-05 14:33:41.419164"
$ws1.Range("F47").Value = "2021-10-05 14:33:41.419167"
$ws1.Range("F48").Value = "2021-10-05 14:33:41.419169"
$ws1.Range("F49").Value = "2021-10-05 14:33:41.419172"
$ws1.Range("F50").Value = "2021-10-05 14:33:41.419174"
$ws1.Range("F51").Value = "2021-10-05 14:33:41.419177"

# --- Add the new "metadata" worksheet, placed after "data" ---
$sheetCount = $wb.Worksheets.Count
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$ws2.Name = "metadata"

# Match the page margins used elsewhere in this workbook (0.75in/1in/0.5in)
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

# Header row (bold, bordered, centered/top-aligned - matches the "data" tab's header style)
$metaHeader = $ws2.Range("B1:G1")
$metaHeader.Font.Bold = $true
$metaHeader.Borders.LineStyle = 1
$metaHeader.HorizontalAlignment = -4108
$metaHeader.VerticalAlignment = -4160

$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Index cell A2 uses the same bold/bordered style as the header and the "data" tab's index column
$metaIndex = $ws2.Range("A2")
$metaIndex.Font.Bold = $true
$metaIndex.Borders.LineStyle = 1
$metaIndex.HorizontalAlignment = -4108
$metaIndex.VerticalAlignment = -4160
$metaIndex.Value = 0

$ws2.Range("B2").Value = "Defects of innate immunity"
$ws2.Range("C2").Value = 231

# Force data_version to stay textual "0.80" rather than being parsed as the number 0.8
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "0.80"

$ws2.Range("E2").Value = "2021-08-25T00:16:24.033145Z"
$ws2.Range("F2").Value = "2021-10-05 14:33:41.415770"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/231/?format=json"

# Keep the "data" tab as the active/selected sheet (matches original activeTab="0")
$ws1.Activate()
